$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) "29.100.72"
Set-TextValue $ws.Cells.Item(2, 5) "  -1.28%  "
Set-TextValue $ws.Cells.Item(3, 4) "1.832.28"
Set-TextValue $ws.Cells.Item(3, 5) "  -1.43%  "
Set-TextValue $ws.Cells.Item(4, 4) "0.9988"
Set-TextValue $ws.Cells.Item(4, 5) "  -0.06%  "
Set-TextValue $ws.Cells.Item(5, 4) "239.06"
Set-TextValue $ws.Cells.Item(5, 5) "  -2.47%  "
Set-TextValue $ws.Cells.Item(6, 4) "0.6632"
Set-TextValue $ws.Cells.Item(6, 5) "  -4.55%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.9997"
Set-TextValue $ws.Cells.Item(8, 4) "0.2947"
Set-TextValue $ws.Cells.Item(8, 5) "  -3.80%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.07330"
Set-TextValue $ws.Cells.Item(9, 5) "  -4.72%  "
Set-TextValue $ws.Cells.Item(10, 4) "22.73"
Set-TextValue $ws.Cells.Item(10, 5) "  -3.81%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.07646"
Set-TextValue $ws.Cells.Item(11, 5) "  -1.57%  "
Set-TextValue $ws.Cells.Item(12, 4) "1.841.17"
Set-TextValue $ws.Cells.Item(12, 5) "  -0.90%  "
Set-TextValue $ws.Cells.Item(13, 4) "5.015"
Set-TextValue $ws.Cells.Item(13, 5) "  -2.66%  "
Set-TextValue $ws.Cells.Item(14, 4) "0.6736"
Set-TextValue $ws.Cells.Item(14, 5) "  -2.72%  "
Set-TextValue $ws.Cells.Item(15, 4) "86.03"
Set-TextValue $ws.Cells.Item(15, 5) "  -5.61%  "
Set-TextValue $ws.Cells.Item(16, 4) "6.126"
Set-TextValue $ws.Cells.Item(16, 5) "  -3.34%  "
Set-TextValue $ws.Cells.Item(17, 4) "29.091.00"
Set-TextValue $ws.Cells.Item(17, 5) "  -1.26%  "
Set-TextValue $ws.Cells.Item(18, 4) "0.000008191"
Set-TextValue $ws.Cells.Item(18, 5) "  -1.29%  "
Set-TextValue $ws.Cells.Item(19, 4) "227.24"
Set-TextValue $ws.Cells.Item(19, 5) "  -4.66%  "
Set-TextValue $ws.Cells.Item(20, 4) "12.47"
Set-TextValue $ws.Cells.Item(20, 5) "  -2.00%  "
Set-TextValue $ws.Cells.Item(21, 4) "0.9994"
Set-TextValue $ws.Cells.Item(21, 5) "  -0.06%  "
Set-TextValue $ws.Cells.Item(22, 4) "7.251"
Set-TextValue $ws.Cells.Item(22, 5) "  -5.11%  "
Set-TextValue $ws.Cells.Item(23, 4) "0.9995"
Set-TextValue $ws.Cells.Item(23, 5) "  -0.06%  "
Set-TextValue $ws.Cells.Item(24, 4) "160.64"
Set-TextValue $ws.Cells.Item(25, 4) "0.1420"
Set-TextValue $ws.Cells.Item(25, 5) "  -4.94%  "
Set-TextValue $ws.Cells.Item(26, 5) "  -2.98%  "
Set-TextValue $ws.Cells.Item(27, 4) "17.93"
Set-TextValue $ws.Cells.Item(27, 5) "  -1.86%  "
Set-TextValue $ws.Cells.Item(28, 4) "1.495"
Set-TextValue $ws.Cells.Item(28, 5) "  -2.46%  "
Set-TextValue $ws.Cells.Item(29, 4) "4.223"
Set-TextValue $ws.Cells.Item(29, 5) "  -0.67%  "
Set-TextValue $ws.Cells.Item(30, 4) "4.100"
Set-TextValue $ws.Cells.Item(30, 5) "  -1.38%  "
Set-TextValue $ws.Cells.Item(31, 4) "1.200"
Set-TextValue $ws.Cells.Item(31, 5) "  -1.26%  "
Set-TextValue $ws.Cells.Item(32, 4) "0.05316"
Set-TextValue $ws.Cells.Item(32, 5) "  +3.95%  "
$ws.Cells.Item(33, 2).Value = "LidoDAOToken"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Cells.Item(33, 4) "1.848"
Set-TextValue $ws.Cells.Item(33, 5) "  -2.04%  "
$ws.Cells.Item(34, 2).Value = "ImmutableX"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Cells.Item(34, 4) "0.7455"
Set-TextValue $ws.Cells.Item(34, 5) "  -3.44%  "
Set-TextValue $ws.Cells.Item(35, 5) "  -1.91%  "
Set-TextValue $ws.Cells.Item(36, 4) "2.678"
Set-TextValue $ws.Cells.Item(36, 5) "  -0.13%  "
Set-TextValue $ws.Cells.Item(37, 4) "1.297.99"
Set-TextValue $ws.Cells.Item(37, 5) "  -2.53%  "
Set-TextValue $ws.Cells.Item(38, 5) "  -3.27%  "
Set-TextValue $ws.Cells.Item(39, 4) "2.704"
Set-TextValue $ws.Cells.Item(39, 5) "  -0.39%  "
Set-TextValue $ws.Cells.Item(40, 4) "0.9209"
Set-TextValue $ws.Cells.Item(40, 5) "  -3.62%  "
Set-TextValue $ws.Cells.Item(41, 4) "6.026"
Set-TextValue $ws.Cells.Item(41, 5) "  +2.85%  "
Set-TextValue $ws.Cells.Item(42, 4) "0.9982"
Set-TextValue $ws.Cells.Item(42, 5) "  -0.22%  "
Set-TextValue $ws.Cells.Item(43, 4) "103.68"
Set-TextValue $ws.Cells.Item(43, 5) "  -2.04%  "
Set-TextValue $ws.Cells.Item(44, 4) "1.983.74"
Set-TextValue $ws.Cells.Item(44, 5) "  -0.83%  "
Set-TextValue $ws.Cells.Item(45, 5) "  -0.93%  "
$ws.Cells.Item(46, 2).Value = "Aave"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Cells.Item(46, 4) "63.42"
Set-TextValue $ws.Cells.Item(46, 5) "  +0.09%  "
Set-TextValue $ws.Cells.Item(47, 4) "1.749"
Set-TextValue $ws.Cells.Item(47, 5) "  -2.03%  "
$ws.Cells.Item(48, 2).Value = "EnergySwap"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Cells.Item(48, 4) "9.210"
Set-TextValue $ws.Cells.Item(48, 5) "  -6.79%  "
$ws.Cells.Item(49, 2).Value = "Cronos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Cells.Item(49, 4) "0.05907"
Set-TextValue $ws.Cells.Item(49, 5) "  -0.38%  "
Set-TextValue $ws.Cells.Item(50, 4) "0.07248"
Set-TextValue $ws.Cells.Item(50, 5) "  +6.45%  "
$ws.Cells.Item(51, 2).Value = "Aptos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Cells.Item(51, 4) "6.817"
Set-TextValue $ws.Cells.Item(51, 5) "  -2.28%  "
